$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 20.447252
$ws.Range("H2").Value = 61.341756
$ws.Range("I2").Value = 0.8699145605694745
$ws.Range("J2").Value = 0.8770588936480435
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 4.472365
$ws.Range("N2").Value = 13.417095
$ws.Range("O2").Value = 0.1840876942178652
$ws.Range("P2").Value = 0.2075728609309428
$ws.Range("Q2").Value = 91.44757419097999
$ws.Range("R2").Value = 823.02816771882
$ws.Range("S2").Value = 0.160140565621782
$ws.Range("T2").Value = 0.1820536237594518

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 20.447252
$ws.Range("H3").Value = 61.341756
$ws.Range("I3").Value = 0.8699145605694745
$ws.Range("J3").Value = 0.8770588936480435
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 3.966196333333333
$ws.Range("N3").Value = 11.898589
$ws.Range("O3").Value = 0.163253208943967
$ws.Range("P3").Value = 0.1840803959256042
$ws.Range("Q3").Value = 81.09781590914265
$ws.Range("R3").Value = 729.8803431822839
$ws.Range("S3").Value = 0.1420163435200477
$ws.Range("T3").Value = 0.1614493483928043

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 20.447252
$ws.Range("H4").Value = 61.341756
$ws.Range("I4").Value = 0.8699145605694745
$ws.Range("J4").Value = 0.8770588936480435
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 4.115150666666667
$ws.Range("N4").Value = 12.345452
$ws.Range("O4").Value = 0.1693843408545093
$ws.Range("P4").Value = 0.1909937129554221
$ws.Range("Q4").Value = 84.14352269930133
$ws.Range("R4").Value = 757.291704293712
$ws.Range("S4").Value = 0.1473499044418006
$ws.Range("T4").Value = 0.1675127345784145

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 20.447252
$ws.Range("H5").Value = 61.341756
$ws.Range("I5").Value = 0.8699145605694745
$ws.Range("J5").Value = 0.8770588936480435
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 3.494784666666666
$ws.Range("N5").Value = 10.484354
$ws.Range("O5").Value = 0.1438493618196675
$ws.Range("P5").Value = 0.1622010841238564
$ws.Range("Q5").Value = 71.45874276506933
$ws.Range("R5").Value = 643.1286848856239
$ws.Range("S5").Value = 0.1251366543755554
$ws.Range("T5").Value = 0.1422599033901827

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 20.447252
$ws.Range("H6").Value = 61.341756
$ws.Range("I6").Value = 0.8699145605694745
$ws.Range("J6").Value = 0.8770588936480435
$ws.Range("K6").Value = 2
$ws.Range("M6").Value = 8.246256
$ws.Range("N6").Value = 16.492512
$ws.Range("O6").Value = 0.3394253941639908
$ws.Range("P6").Value = 0.2551519460641745
$ws.Range("Q6").Value = 168.613274488512
$ws.Range("R6").Value = 1011.679646931072
$ws.Range("S6").Value = 0.2952710926102887
$ws.Range("T6").Value = 0.2237832835271901

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 2.483247333333333
$ws.Range("H7").Value = 7.449742
$ws.Range("I7").Value = 0.105648084777455
$ws.Range("J7").Value = 0.1065157390747562
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 4.472365
$ws.Range("N7").Value = 13.417095
$ws.Range("O7").Value = 0.1840876942178652
$ws.Range("P7").Value = 0.2075728609309428
$ws.Range("Q7").Value = 11.10598845994333
$ws.Range("R7").Value = 99.95389613949
$ws.Range("S7").Value = 0.01944851232521524
$ws.Range("T7").Value = 0.02210977669392096

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 2.483247333333333
$ws.Range("H8").Value = 7.449742
$ws.Range("I8").Value = 0.105648084777455
$ws.Range("J8").Value = 0.1065157390747562
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 3.966196333333333
$ws.Range("N8").Value = 11.898589
$ws.Range("O8").Value = 0.163253208943967
$ws.Range("P8").Value = 0.1840803959256042
$ws.Range("Q8").Value = 9.849046468226444
$ws.Range("R8").Value = 88.641418214038
$ws.Range("S8").Value = 0.0172473888587038
$ws.Range("T8").Value = 0.01960745942118948

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 2.483247333333333
$ws.Range("H9").Value = 7.449742
$ws.Range("I9").Value = 0.105648084777455
$ws.Range("J9").Value = 0.1065157390747562
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 4.115150666666667
$ws.Range("N9").Value = 12.345452
$ws.Range("O9").Value = 0.1693843408545093
$ws.Range("P9").Value = 0.1909937129554221
$ws.Range("Q9").Value = 10.21893691926489
$ws.Range("R9").Value = 91.97043227338399
$ws.Range("S9").Value = 0.01789513120257054
$ws.Range("T9").Value = 0.02034383649407863

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 2.483247333333333
$ws.Range("H10").Value = 7.449742
$ws.Range("I10").Value = 0.105648084777455
$ws.Range("J10").Value = 0.1065157390747562
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 3.494784666666666
$ws.Range("N10").Value = 10.484354
$ws.Range("O10").Value = 0.1438493618196675
$ws.Range("P10").Value = 0.1622010841238564
$ws.Range("Q10").Value = 8.678414704074221
$ws.Range("R10").Value = 78.105732336668
$ws.Range("S10").Value = 0.01519740957270703
$ws.Range("T10").Value = 0.01727696835417928

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 2.483247333333333
$ws.Range("H11").Value = 7.449742
$ws.Range("I11").Value = 0.105648084777455
$ws.Range("J11").Value = 0.1065157390747562
$ws.Range("K11").Value = 2
$ws.Range("M11").Value = 8.246256
$ws.Range("N11").Value = 16.492512
$ws.Range("O11").Value = 0.3394253941639908
$ws.Range("P11").Value = 0.2551519460641745
$ws.Range("Q11").Value = 20.477493221984
$ws.Range("R11").Value = 122.864959331904
$ws.Range("S11").Value = 0.03585964281825838
$ws.Range("T11").Value = 0.02717769811138789

$ws.Range("E12").Value = 2
$ws.Range("G12").Value = 0.5743975
$ws.Range("H12").Value = 1.148795
$ws.Range("I12").Value = 0.02443735465307048
$ws.Range("J12").Value = 0.01642536727720028
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 4.472365
$ws.Range("N12").Value = 13.417095
$ws.Range("O12").Value = 0.1840876942178652
$ws.Range("P12").Value = 0.2075728609309428
$ws.Range("Q12").Value = 2.5689152750875
$ws.Range("R12").Value = 15.413491650525
$ws.Range("S12").Value = 0.004498616270867966
$ws.Range("T12").Value = 0.003409460477569952

$ws.Range("E13").Value = 2
$ws.Range("G13").Value = 0.5743975
$ws.Range("H13").Value = 1.148795
$ws.Range("I13").Value = 0.02443735465307048
$ws.Range("J13").Value = 0.01642536727720028
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 3.966196333333333
$ws.Range("N13").Value = 11.898589
$ws.Range("O13").Value = 0.163253208943967
$ws.Range("P13").Value = 0.1840803959256042
$ws.Range("Q13").Value = 2.278173258375833
$ws.Range("R13").Value = 13.669039550255
$ws.Range("S13").Value = 0.00398947656521554
$ws.Range("T13").Value = 0.003023588111610492

$ws.Range("E14").Value = 2
$ws.Range("G14").Value = 0.5743975
$ws.Range("H14").Value = 1.148795
$ws.Range("I14").Value = 0.02443735465307048
$ws.Range("J14").Value = 0.01642536727720028
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 4.115150666666667
$ws.Range("N14").Value = 12.345452
$ws.Range("O14").Value = 0.1693843408545093
$ws.Range("P14").Value = 0.1909937129554221
$ws.Range("Q14").Value = 2.363732255056667
$ws.Range("R14").Value = 14.18239353034
$ws.Range("S14").Value = 0.004139305210138221
$ws.Range("T14").Value = 0.003137141882928974

$ws.Range("E15").Value = 2
$ws.Range("G15").Value = 0.5743975
$ws.Range("H15").Value = 1.148795
$ws.Range("I15").Value = 0.02443735465307048
$ws.Range("J15").Value = 0.01642536727720028
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 3.494784666666666
$ws.Range("N15").Value = 10.484354
$ws.Range("O15").Value = 0.1438493618196675
$ws.Range("P15").Value = 0.1622010841238564
$ws.Range("Q15").Value = 2.007395575571667
$ws.Range("R15").Value = 12.04437345343
$ws.Range("S15").Value = 0.003515297871405072
$ws.Range("T15").Value = 0.002664212379494402

$ws.Range("E16").Value = 2
$ws.Range("G16").Value = 0.5743975
$ws.Range("H16").Value = 1.148795
$ws.Range("I16").Value = 0.02443735465307048
$ws.Range("J16").Value = 0.01642536727720028
$ws.Range("K16").Value = 2
$ws.Range("M16").Value = 8.246256
$ws.Range("N16").Value = 16.492512
$ws.Range("O16").Value = 0.3394253941639908
$ws.Range("P16").Value = 0.2551519460641745
$ws.Range("Q16").Value = 4.736628830760001
$ws.Range("R16").Value = 18.94651532304
$ws.Range("S16").Value = 0.008294658735443685
$ws.Range("T16").Value = 0.004190964425596463
